$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-NumCell($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value2 = [double]$value
}

function Set-TextCell($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $styleSrc = $ws.Cells.Item($row, 1)
    $styleSrc.Copy()
    $cell.PasteSpecial(-4122)
}

# --- Header text updates ---

# "Volume 31   Number  37" -> "...38"
$volCell = $ws.Range("A8")
$volText = $volCell.Value2
$idx = $volText.LastIndexOf("37")
$volCell.Characters($idx + 1, 2).Text = "38"

# "Report Covering the Week  9/9/2024  Through  9/15/2024"
#   -> "...9/16/2024  Through  9/22/2024"
$dateCell = $ws.Range("C9")
$dateText = $dateCell.Value2
$i1 = $dateText.IndexOf("9/9/2024")
$dateCell.Characters($i1 + 1, 8).Text = "9/16/2024"
$dateText2 = $dateCell.Value2
$i2 = $dateText2.IndexOf("9/15/2024")
$dateCell.Characters($i2 + 1, 9).Text = "9/22/2024"

$excel.CutCopyMode = $false

# --- Crime statistics table updates (rows 14-30) ---

Set-TextCell 14 4 "0"
Set-TextCell 14 5 "***.*"
Set-NumCell 14 12 16.666666666666
Set-TextCell 15 3 "0"
Set-NumCell 15 6 2
Set-TextCell 15 7 "0"
Set-TextCell 15 8 "***.*"
Set-NumCell 15 14 -62.068965517241
Set-NumCell 16 3 2
Set-NumCell 16 4 11
Set-NumCell 16 5 -81.818181818181
Set-NumCell 16 6 9
Set-NumCell 16 7 32
Set-NumCell 16 8 -71.875
Set-NumCell 16 9 183
Set-NumCell 16 10 203
Set-NumCell 16 11 -9.852216748768
Set-NumCell 16 12 3.977272727272
Set-NumCell 16 13 9.580838323353
Set-NumCell 16 14 -71.271585557299
Set-NumCell 17 3 10
Set-NumCell 17 5 25
Set-NumCell 17 6 36
Set-NumCell 17 7 25
Set-NumCell 17 8 44
Set-NumCell 17 9 362
Set-NumCell 17 10 288
Set-NumCell 17 11 25.694444444444
Set-NumCell 17 12 4.322766570605
Set-NumCell 17 13 110.46511627907
Set-NumCell 17 14 -28.031809145129
Set-NumCell 18 3 4
Set-NumCell 18 4 3
Set-NumCell 18 5 33.333333333333
Set-NumCell 18 6 12
Set-NumCell 18 7 10
Set-NumCell 18 8 20
Set-NumCell 18 9 76
Set-NumCell 18 10 88
Set-NumCell 18 11 -13.636363636363
Set-NumCell 18 12 -51.282051282051
Set-NumCell 18 13 -6.172839506172
Set-NumCell 18 14 -80.25974025974
Set-NumCell 19 3 12
Set-NumCell 19 4 15
Set-NumCell 19 5 -20
Set-NumCell 19 6 38
Set-NumCell 19 7 115
Set-NumCell 19 8 -66.95652173913
Set-NumCell 19 9 378
Set-NumCell 19 10 429
Set-NumCell 19 11 -11.888111888111
Set-NumCell 19 12 -19.574468085106
Set-NumCell 19 13 136.25
Set-NumCell 19 14 29.010238907849
Set-TextCell 20 3 "0"
Set-NumCell 20 4 3
Set-NumCell 20 5 -100
Set-NumCell 20 7 12
Set-NumCell 20 8 -41.666666666666
Set-NumCell 20 10 66
Set-NumCell 20 11 27.272727272727
Set-NumCell 20 13 55.555555555555
Set-NumCell 20 14 -65.714285714285
Set-NumCell 21 3 28
Set-NumCell 21 4 40
Set-NumCell 21 5 -30
Set-NumCell 21 6 104
Set-NumCell 21 7 195
Set-NumCell 21 8 -46.666666666666
Set-NumCell 21 9 1101
Set-NumCell 21 10 1089
Set-NumCell 21 11 1.101928374655
Set-NumCell 21 12 -10.268948655256
Set-NumCell 21 13 68.091603053435
Set-NumCell 21 14 -48.066037735849
Set-TextCell 22 7 "0"
Set-TextCell 22 8 "***.*"
Set-NumCell 23 3 3
Set-NumCell 23 4 3
Set-NumCell 23 7 10
Set-NumCell 23 8 -20
Set-NumCell 23 9 121
Set-NumCell 23 10 120
Set-NumCell 23 11 0.833333333333
Set-NumCell 23 12 -9.701492537313
Set-NumCell 23 13 45.78313253012
Set-NumCell 24 3 22
Set-NumCell 24 4 21
Set-NumCell 24 5 4.761904761904
Set-NumCell 24 6 66
Set-NumCell 24 7 83
Set-NumCell 24 8 -20.481927710843
Set-NumCell 24 9 734
Set-NumCell 24 10 917
Set-NumCell 24 11 -19.956379498364
Set-NumCell 24 12 -22.410147991543
Set-NumCell 24 13 5.308464849354
Set-NumCell 25 3 2
Set-NumCell 25 4 7
Set-NumCell 25 5 -71.428571428571
Set-NumCell 25 6 11
Set-NumCell 25 7 40
Set-NumCell 25 8 -72.5
Set-NumCell 25 9 145
Set-NumCell 25 10 443
Set-NumCell 25 11 -67.26862302483
Set-NumCell 25 12 -73.247232472324
Set-NumCell 26 3 18
Set-NumCell 26 4 11
Set-NumCell 26 5 63.636363636363
Set-NumCell 26 6 65
Set-NumCell 26 7 47
Set-NumCell 26 8 38.297872340425
Set-NumCell 26 9 572
Set-NumCell 26 10 377
Set-NumCell 26 11 51.724137931034
Set-NumCell 26 12 36.515513126491
Set-NumCell 26 13 56.284153005464
Set-TextCell 27 3 "0"
Set-TextCell 27 4 "0"
Set-TextCell 27 5 "***.*"
Set-NumCell 27 6 2
Set-NumCell 27 7 1
Set-NumCell 28 3 2
Set-NumCell 28 4 4
Set-NumCell 28 5 -50
Set-NumCell 28 6 8
Set-NumCell 28 7 10
Set-NumCell 28 8 -20
Set-NumCell 28 9 48
Set-NumCell 28 10 46
Set-NumCell 28 11 4.347826086956
Set-NumCell 28 12 -11.111111111111
Set-TextCell 29 3 "0"
Set-NumCell 29 5 -100
Set-TextCell 29 6 "0"
Set-NumCell 29 7 4
Set-NumCell 29 8 -100
Set-NumCell 29 9 12
Set-NumCell 29 10 27
Set-NumCell 29 11 -55.555555555555
Set-NumCell 29 12 -66.666666666666
Set-NumCell 29 13 -58.620689655172
Set-NumCell 29 14 -78.571428571428
Set-TextCell 30 3 "0"
Set-NumCell 30 5 -100
Set-TextCell 30 6 "0"
Set-NumCell 30 7 3
Set-NumCell 30 8 -100
Set-NumCell 30 9 9
Set-NumCell 30 10 25
Set-NumCell 30 11 -64
Set-NumCell 30 12 -65.384615384615
Set-NumCell 30 13 -62.5
Set-NumCell 30 14 -83.018867924528

$excel.CutCopyMode = $false
